$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.347.34"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.925.26"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.91"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.02"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.93"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.48"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "3.411.11"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "61.272.88"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "2.928.36"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.11"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.675"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.89"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.76"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.64"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "0.0₃0879"
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.63"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.98"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.00"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.44"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.281"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0345"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "2.697.07"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "366.59"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.57"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.54"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.124"
$ws.Range("E51").Value = "  -1.02%  "
